$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Matching")

# Populate the two new matching pairs (Hot/Fire and Cold/Ice, both under
# the "temperature" category) in the order Excel needs so the shared
# strings table comes out Hot, Cold, Fire, Ice, temperature.
$ws.Range("A5").Value = "Hot"
$ws.Range("A6").Value = "Cold"
$ws.Range("B5").Value = "Fire"
$ws.Range("B6").Value = "Ice"
$ws.Range("C5").Value = "temperature"
$ws.Range("C6").Value = "temperature"

# Column C has no inherited cell style, so match the top-aligned style
# already used by columns A/B in this block.
$ws.Range("C5:C6").VerticalAlignment = -4160

# Make "Matching" the active sheet/tab, with C6 as the selected cell -
# mirrors the author re-opening the workbook on that sheet.
$ws.Select()
$ws.Range("C6").Select()
